$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values with newly repulled/recalculated data.
$ws.Range("F2").Value = -14
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("F12").Value = -5
$ws.Range("F13").Value = -4
$ws.Range("F14").Value = -9
$ws.Range("F16").Value = -4
